$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DescriptivesAndQuestionnaires")

# Insert a new column before the old "Counterbalancing order*" column (D).
# This shifts old D->E, E->F, F->G, G->H, H->I, I->J, K->L, matching the diff.
$ws.Columns("D:D").Insert()

# Rename old "Age" header (column C) to "Date of Birth", and label the new
# column D "Date of Experiment".
$ws.Range("C1").Value = "Date of Birth"
$ws.Range("D1").Value = "Date of Experiment"

# Both new data columns hold Excel date serials formatted as dates.
$ws.Range("C2:D33").NumberFormat = "m/d/yy"

$dates = @{
    2 = @(34074, 42944)
    3 = @(34660, 42948)
    4 = @(32404, 42948)
    5 = @(31070, 42949)
    6 = @(34336, 42950)
    7 = @(33405, 42950)
    8 = @(30727, 42950)
    9 = @(34435, 42951)
    10 = @(33761, 42951)
    11 = @(33382, 42951)
    12 = @(34307, 42951)
    13 = @(33535, 42954)
    14 = @(33990, 42954)
    15 = @(32476, 42955)
    16 = @(35338, 42955)
    17 = @(35301, 42955)
    18 = @(33099, 42956)
    19 = @(34154, 42956)
    20 = @(33782, 42956)
    21 = @(33956, 42956)
    22 = @(29148, 42956)
    23 = @(32739, 42957)
    24 = @(33251, 42958)
    25 = @(33244, 42958)
    26 = @(32700, 42968)
    27 = @(31989, 42968)
    28 = @(31655, 42969)
    29 = @(32837, 42975)
    30 = @(31726, 42976)
    31 = @(33313, 42976)
    32 = @(33478, 42977)
    33 = @(33844, 42977)
}

foreach ($r in $dates.Keys) {
    $pair = $dates[$r]
    $ws.Cells.Item($r, 3).Value = $pair[0]
    $ws.Cells.Item($r, 4).Value = $pair[1]
}
